$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3492.3044
$ws.Range("I15").Value = 3492.3044
$ws.Range("K15").Value = 10476.9132
$ws.Range("M15").Value = -10307.9132

$ws.Range("H76").Value = 47623610
$ws.Range("I76").Value = 3414.4167
$ws.Range("J76").Value = 111117200
$ws.Range("K76").Value = 3414.4167
$ws.Range("L76").Value = 111117200
$ws.Range("M76").Value = -3099.4167
$ws.Range("N76").Value = -111117830

$ws.Range("H79").Value = 47623610
$ws.Range("I79").Value = 3414.4167
$ws.Range("J79").Value = 111117200
$ws.Range("K79").Value = 3414.4167
$ws.Range("L79").Value = 111117200
$ws.Range("M79").Value = -2322.4167
$ws.Range("N79").Value = -111119384

$ws.Range("H112").Value = 1977.7428
$ws.Range("J112").Value = 2003.5588
$ws.Range("L112").Value = 6010.6764
$ws.Range("N112").Value = -8226.6764

$ws.Range("H113").Value = 3690.6
$ws.Range("I113").Value = 2900
$ws.Range("J113").Value = 3778.4443
$ws.Range("K113").Value = 2900
$ws.Range("L113").Value = 3778.4443
$ws.Range("M113").Value = 354
$ws.Range("N113").Value = -10286.4443

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3552.077
$ws.Range("I32").Value = 3275.6956
$ws.Range("J32").Value = 5671
$ws.Range("K32").Value = 3275.6956
$ws.Range("L32").Value = 5671
$ws.Range("M32").Value = -2988.6956
$ws.Range("N32").Value = -6245

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3934.28
$ws.Range("I86").Value = 4408.9443
$ws.Range("J86").Value = 2713.7144
$ws.Range("K86").Value = 4408.9443
$ws.Range("L86").Value = 2713.7144
$ws.Range("M86").Value = -3285.9443
$ws.Range("N86").Value = -4959.7144

$ws.Range("H89").Value = 3934.28
$ws.Range("I89").Value = 4408.9443
$ws.Range("J89").Value = 2713.7144
$ws.Range("K89").Value = 22044.7215
$ws.Range("L89").Value = 13568.572
$ws.Range("M89").Value = -16428.7215
$ws.Range("N89").Value = -24800.572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1184.1333
$ws.Range("I31").Value = 1054.4286
$ws.Range("K31").Value = 1054.4286
$ws.Range("M31").Value = -759.4286

$ws.Range("H34").Value = 1184.1333
$ws.Range("I34").Value = 1054.4286
$ws.Range("K34").Value = 1054.4286
$ws.Range("M34").Value = -852.4286

$ws.Range("H99").Value = 1962
$ws.Range("I99").Value = 1899.3334
$ws.Range("J99").Value = 2150
$ws.Range("K99").Value = 1899.3334
$ws.Range("L99").Value = 2150
$ws.Range("M99").Value = -401.3334
$ws.Range("N99").Value = -5146

$ws.Range("H107").Value = 534
$ws.Range("I107").Value = 353.88235
$ws.Range("J107").Value = 874.2222
$ws.Range("K107").Value = 353.88235
$ws.Range("L107").Value = 874.2222
$ws.Range("M107").Value = 1566.11765
$ws.Range("N107").Value = -4714.2222

$ws.Range("H126").Value = 1962
$ws.Range("I126").Value = 1899.3334
$ws.Range("J126").Value = 2150
$ws.Range("K126").Value = 5698.0002
$ws.Range("L126").Value = 6450
$ws.Range("M126").Value = -3228.0002
$ws.Range("N126").Value = -11390

$ws.Range("H132").Value = 5571.423
$ws.Range("I132").Value = 6609.3887
$ws.Range("J132").Value = 3236
$ws.Range("K132").Value = 19828.1661
$ws.Range("L132").Value = 9708
$ws.Range("M132").Value = -17298.1661
$ws.Range("N132").Value = -14768

$ws.Range("H134").Value = 2123.6956
$ws.Range("I134").Value = 2340.4666
$ws.Range("J134").Value = 1717.25
$ws.Range("K134").Value = 7021.399800000001
$ws.Range("L134").Value = 5151.75
$ws.Range("M134").Value = -4486.399800000001
$ws.Range("N134").Value = -10221.75

$ws.Range("H135").Value = 62223.332
$ws.Range("J135").Value = 62223.332
$ws.Range("L135").Value = 62223.332
$ws.Range("N135").Value = -72363.33199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 780421.75
$ws.Range("J4").Value = 900598.9399999999
$ws.Range("L4").Value = 2701796.82
$ws.Range("N4").Value = -2702020.82

$ws.Range("H34").Value = 12501405
$ws.Range("I34").Value = 412.33334
$ws.Range("J34").Value = 20002000
$ws.Range("K34").Value = 1237.00002
$ws.Range("L34").Value = 60006000
$ws.Range("M34").Value = -1153.00002
$ws.Range("N34").Value = -60006168

$ws.Range("H46").Value = 3003
$ws.Range("I46").Value = 3003
$ws.Range("K46").Value = 9009
$ws.Range("M46").Value = -8918

$ws.Range("H96").Value = 8142.7144
$ws.Range("J96").Value = 8142.7144
$ws.Range("L96").Value = 24428.1432
$ws.Range("N96").Value = -28546.1432

$ws.Range("H113").Value = 603.5135
$ws.Range("I113").Value = 400
$ws.Range("J113").Value = 659.65515
$ws.Range("K113").Value = 1200
$ws.Range("L113").Value = 1978.96545
$ws.Range("M113").Value = 970
$ws.Range("N113").Value = -6318.96545

$ws.Range("H131").Value = 14085634
$ws.Range("I131").Value = 100000400
$ws.Range("J131").Value = 1245.7705
$ws.Range("K131").Value = 300001200
$ws.Range("L131").Value = 3737.3115
$ws.Range("M131").Value = -299996160
$ws.Range("N131").Value = -13817.3115

$ws.Range("H138").Value = 4573.8
$ws.Range("I138").Value = 5017.25
$ws.Range("K138").Value = 15051.75
$ws.Range("M138").Value = -9911.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2361.4644
$ws.Range("I102").Value = 2895.842
$ws.Range("J102").Value = 1233.3334
$ws.Range("K102").Value = 2895.842
$ws.Range("L102").Value = 1233.3334
$ws.Range("M102").Value = -1273.842
$ws.Range("N102").Value = -4477.3334

$ws.Range("H122").Value = 1240.909
$ws.Range("I122").Value = 1318.75
$ws.Range("K122").Value = 3956.25
$ws.Range("M122").Value = -1506.25

$ws.Range("H126").Value = 2155.8235
$ws.Range("I126").Value = 2215
$ws.Range("J126").Value = 2089.25
$ws.Range("K126").Value = 6645
$ws.Range("L126").Value = 6267.75
$ws.Range("M126").Value = -4175
$ws.Range("N126").Value = -11207.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2848.3
$ws.Range("I40").Value = 2782.3333
$ws.Range("J40").Value = 2947.25
$ws.Range("K40").Value = 2782.3333
$ws.Range("L40").Value = 2947.25
$ws.Range("M40").Value = -2646.3333
$ws.Range("N40").Value = -3219.25

$ws.Range("H106").Value = 30712.715
$ws.Range("J106").Value = 30712.715
$ws.Range("L106").Value = 30712.715
$ws.Range("N106").Value = -33236.715

$ws.Range("H122").Value = 13160468
$ws.Range("I122").Value = 27780228
$ws.Range("J122").Value = 2683.9
$ws.Range("K122").Value = 83340684
$ws.Range("L122").Value = 8051.700000000001
$ws.Range("M122").Value = -83338234
$ws.Range("N122").Value = -12951.7

$ws.Range("H132").Value = 20860.424
$ws.Range("I132").Value = 1360.3704
$ws.Range("J132").Value = 41920.48
$ws.Range("K132").Value = 4081.1112
$ws.Range("L132").Value = 125761.44
$ws.Range("M132").Value = -1551.1112
$ws.Range("N132").Value = -130821.44

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10066.25
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 10066.25
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 10066.25
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -10846.25

$ws.Range("H104").Value = 27123
$ws.Range("J104").Value = 27123
$ws.Range("L104").Value = 27123
$ws.Range("N104").Value = -34111

$ws.Range("H122").Value = 8389092
$ws.Range("I122").Value = 10002052
$ws.Range("J122").Value = 1700
$ws.Range("K122").Value = 30006156
$ws.Range("L122").Value = 5100
$ws.Range("M122").Value = -30003706
$ws.Range("N122").Value = -10000

$ws.Range("H132").Value = 2584.9062
$ws.Range("I132").Value = 2211.087
$ws.Range("J132").Value = 3540.2222
$ws.Range("K132").Value = 6633.261
$ws.Range("L132").Value = 10620.6666
$ws.Range("M132").Value = -4103.261
$ws.Range("N132").Value = -15680.6666
